$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.890.90'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -2.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.898.83'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -4.82%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.58'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4591'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -2.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3807'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -3.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.61'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -2.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07736'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -3.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9803'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -2.29%  '

$ws.Range("E12").Value = '  -3.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.924.12'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -3.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.959'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -4.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.671'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -3.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07067'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -1.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '84.02'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -5.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009537'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -4.83%  '

$ws.Range("E20").Value = '  -4.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.006'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.846.90'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("E23").Value = '  -4.11%  '

$ws.Range("E24").Value = '  -3.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.152.74'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -4.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.102'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -0.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.33'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -0.55%  '

$ws.Range("E28").Value = '  -2.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.585'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -7.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.64'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -2.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.837'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -5.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09268'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -2.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8615'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -6.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.099'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -3.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.252'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -7.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.024'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -5.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05693'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.146'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -2.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.006'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02035'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -4.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.477'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -5.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5511'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -4.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1755'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -4.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.321'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -5.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.722'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5194'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -3.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.24'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -7.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.101'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -4.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06835'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -1.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.56'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -2.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.772'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -5.73%  '
